$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the newly highlighted "auto-consequence" cells in the matrix (yellow fill) ---
$ws.Range("H5").Value = "SBS"
$ws.Range("H5").Interior.Color = 65535

$ws.Range("H6").Value = "SBS"
$ws.Range("H6").Interior.Color = 65535

$ws.Range("H8").Value = "FBS"
$ws.Range("H8").Interior.Color = 65535

$ws.Range("H9").Value = "FBS"
$ws.Range("H9").Interior.Color = 65535

$ws.Range("I10").Value = "SAF"
$ws.Range("I10").Interior.Color = 65535

$ws.Range("I11").Value = "SAF"
$ws.Range("I11").Interior.Color = 65535

$ws.Range("I13").Value = "FAF"
$ws.Range("I13").Interior.Color = 65535

$ws.Range("I14").Value = "FAF"
$ws.Range("I14").Interior.Color = 65535

# --- Add the new "Auto-consequences" legend block below the matrix ---
$ws.Range("B21").Value = "SBS > SBF"
$ws.Range("B21").Interior.Color = 65535

$ws.Range("B24").Value = "FAF > FAS"
$ws.Range("B24").Interior.Color = 65535

$ws.Range("B20").Value = "Auto-consequences:"
$ws.Range("B20").Interior.Color = 65535

$ws.Range("B23").Value = "SAF > FAF, FAS, SAS"
$ws.Range("B23").Interior.Color = 65535

$ws.Range("B22").Value = "FBS > SBS, SBF, FBF"
$ws.Range("B22").Interior.Color = 65535

# --- Update the view: select the newly added legend block ---
$ws.Range("B21:B24").Select()
